# Update the "Förändrad" (Changed) date column (C) for all data rows.
# Every data row (2 through 295) currently holds the Excel serial date
# 45179 (2023-09-10) in column C; bump it to 45180 (2023-09-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 295
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45179) {
        $cell.Value2 = 45180
    }
}
